# WebForm User Assignment execution
#
# The WebForm automation run assigned phone numbers (PN_Value, column F)
# to rows 2-10 and reset the Match2UserPos / UnMatchUserPos counters
# (AN2 / AO2) back to "0".
#
# Values are written with a leading apostrophe so Excel keeps them as
# text (these are phone numbers, not numeric quantities) instead of
# silently re-typing them as numbers; re-applying the "Normal" style
# afterwards clears the resulting quote-prefix flag so the cells keep
# their original (default) formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = "'9840057845"
$ws.Range("F3").Value  = "'9840031971"
$ws.Range("F4").Value  = "'9840046938"
$ws.Range("F5").Value  = "'9840045705"
$ws.Range("F6").Value  = "'9840037612"
$ws.Range("F7").Value  = "'9840057511"
$ws.Range("F8").Value  = "'9840059515"
$ws.Range("F9").Value  = "'9840012168"
$ws.Range("F10").Value = "'9840018722"
$ws.Range("F2:F10").Style = "Normal"

$ws.Range("AN2").Value = "'0"
$ws.Range("AO2").Value = "'0"
$ws.Range("AN2:AO2").Style = "Normal"
